$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename recipe_code values in column A (REC_VACHE_BRASSE_* -> REC_PROC_BASE_VACHE_BRASSE_*)
$ws.Range("A2").Value = "REC_PROC_BASE_VACHE_BRASSE_NATURE"
$ws.Range("A3").Value = "REC_PROC_BASE_VACHE_BRASSE_NATURE"
$ws.Range("A4").Value = "REC_PROC_BASE_VACHE_BRASSE_NATURE"
$ws.Range("A5").Value = "REC_PROC_BASE_VACHE_BRASSE_SUCRE"
$ws.Range("A6").Value = "REC_PROC_BASE_VACHE_BRASSE_SUCRE"

# Add new row 7 (new recipe entry REC_COND_1025700)
$ws.Range("A7").Value = "REC_COND_1025700"
$ws.Range("B7").Value = "BASE_VACHE_BRASSE_NATURE"
$ws.Range("C7").Value = 1000
$ws.Range("D7").Value = "L"

# E7 should hold the text "True" (like E2:E6) rather than a boolean -
# copy an existing "True" text cell so the type/shared-string is preserved.
$ws.Range("E2").Copy()
$ws.Range("E7").PasteSpecial()
$excel.CutCopyMode = 0

# Widen column A to fit the longer recipe_code values
$ws.Columns.Item(1).ColumnWidth = 47.333333333333336

# Move the active selection, matching the saved view state
$ws.Range("E9").Select() | Out-Null
